$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: hours bumped from 2 to 4, with the slightly tighter row height
#     Excel ends up applying after re-entering the row ---
$ws.Range("B4").Value = 4
$ws.Rows.Item(4).RowHeight = 16.5

# --- New log entries in rows 8-11 ---
$ws.Range("A8").Value = "Implemented Shader Toggle"
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = 43988
$ws.Range("D8").Value = "Implemented Shader toggle to switch between custom shader & default URP"

$ws.Range("A9").Value = "Animated Small menu & added extra Walls"
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = 43988
$ws.Range("D9").Value = "Added menu animations & added more walls variants"

$ws.Range("A10").Value = "Made Documentation"
$ws.Range("B10").Value = 1
$ws.Range("C10").Value = 43988
$ws.Range("D10").Value = "Used doxygen for documentation page and made UML Diagrams"

$ws.Range("A11").Value = "Bug Fixes"
$ws.Range("B11").Value = 1
$ws.Range("C11").Value = 43988
$ws.Range("D11").Value = "fixed some background calculation bugs"

# --- Scroll the view down a bit and leave the selection where the author left it ---
$ws.Activate()
$ws.Range("E10").Select()
